# Auto-generated edit script: updates phantom-profit columns (H-N)
# across multiple sheets to match the scheduled-runner recompute.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets("ALC")
$ws.Range("H9").Value = 929.93335
$ws.Range("I9").Value = 987.4167
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 987.4167
$ws.Range("L9").Value = 700
$ws.Range("M9").Value = -818.4167
$ws.Range("N9").Value = -1038
$ws.Range("H40").Value = 1870
$ws.Range("I40").Value = 1863.75
$ws.Range("J40").Value = 1882.5
$ws.Range("K40").Value = 1863.75
$ws.Range("L40").Value = 1882.5
$ws.Range("M40").Value = -1688.75
$ws.Range("N40").Value = -2232.5
$ws.Range("H41").Value = 303.3846
$ws.Range("J41").Value = 308
$ws.Range("L41").Value = 308
$ws.Range("N41").Value = -1188
$ws.Range("H53").Value = 243.28572
$ws.Range("I53").Value = 259.72726
$ws.Range("K53").Value = 259.72726
$ws.Range("M53").Value = 377.27274
$ws.Range("H80").Value = 1393.1875
$ws.Range("I80").Value = 749.2
$ws.Range("J80").Value = 1685.909
$ws.Range("K80").Value = 2247.6
$ws.Range("L80").Value = 5057.727000000001
$ws.Range("M80").Value = -1249.6
$ws.Range("N80").Value = -7053.727000000001
$ws.Range("H83").Value = 1393.1875
$ws.Range("I83").Value = 749.2
$ws.Range("J83").Value = 1685.909
$ws.Range("K83").Value = 6742.8
$ws.Range("L83").Value = 15173.181
$ws.Range("M83").Value = -1750.8
$ws.Range("N83").Value = -25157.181
$ws.Range("H106").Value = 5342.643
$ws.Range("I106").Value = 5676.6924
$ws.Range("K106").Value = 5676.6924
$ws.Range("M106").Value = -5045.6924
$ws.Range("H113").Value = 3633
$ws.Range("I113").Value = 3633
$ws.Range("K113").Value = 3633
$ws.Range("M113").Value = -379
$ws.Range("H125").Value = 1597.0667
$ws.Range("I125").Value = 2067.3333
$ws.Range("J125").Value = 1479.5
$ws.Range("K125").Value = 18605.9997
$ws.Range("L125").Value = 13315.5
$ws.Range("M125").Value = -16145.9997
$ws.Range("N125").Value = -18235.5
$ws.Range("H137").Value = 3046.8333
$ws.Range("I137").Value = 1284.4
$ws.Range("J137").Value = 5249.875
$ws.Range("K137").Value = 3853.2
$ws.Range("L137").Value = 15749.625
$ws.Range("M137").Value = -1303.2
$ws.Range("N137").Value = -20849.625
$ws.Range("H138").Value = 1869.68
$ws.Range("J138").Value = 3299
$ws.Range("L138").Value = 9897
$ws.Range("N138").Value = -20177

$ws = $wb.Sheets("ARM")
$ws.Range("H2").Value = 1046
$ws.Range("I2").Value = 910
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 910
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -797
$ws.Range("N2").Value = -1476
$ws.Range("H5").Value = 557.125
$ws.Range("I5").Value = 595.4286
$ws.Range("K5").Value = 595.4286
$ws.Range("M5").Value = -483.4286
$ws.Range("H43").Value = 31000
$ws.Range("J43").Value = 31000
$ws.Range("L43").Value = 31000
$ws.Range("N43").Value = -31626
$ws.Range("H61").Value = 6252.483
$ws.Range("I61").Value = 5557.2915
$ws.Range("K61").Value = 5557.2915
$ws.Range("M61").Value = -5345.2915
$ws.Range("H116").Value = 1046
$ws.Range("I116").Value = 910
$ws.Range("J116").Value = 1250
$ws.Range("K116").Value = 910
$ws.Range("L116").Value = 1250
$ws.Range("M116").Value = 1384
$ws.Range("N116").Value = -5838
$ws.Range("H136").Value = 6252.483
$ws.Range("I136").Value = 5557.2915
$ws.Range("K136").Value = 16671.8745
$ws.Range("M136").Value = -14121.8745

$ws = $wb.Sheets("BSM")
$ws.Range("H3").Value = 1046
$ws.Range("I3").Value = 910
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 910
$ws.Range("L3").Value = 1250
$ws.Range("M3").Value = -796
$ws.Range("N3").Value = -1478
$ws.Range("H4").Value = 557.125
$ws.Range("I4").Value = 595.4286
$ws.Range("K4").Value = 595.4286
$ws.Range("M4").Value = -480.4286
$ws.Range("H99").Value = 597.0909
$ws.Range("I99").Value = 575.7
$ws.Range("K99").Value = 575.7
$ws.Range("M99").Value = 922.3
$ws.Range("H105").Value = 4765.8
$ws.Range("I105").Value = 4739.778
$ws.Range("K105").Value = 4739.778
$ws.Range("M105").Value = -2992.778
$ws.Range("H112").Value = 48015
$ws.Range("I112").Value = 48015
$ws.Range("K112").Value = 48015
$ws.Range("M112").Value = -46538
$ws.Range("H134").Value = 5062.304
$ws.Range("I134").Value = 4887.409
$ws.Range("K134").Value = 14662.227
$ws.Range("M134").Value = -12127.227

$ws = $wb.Sheets("CRP")
$ws.Range("H7").Value = 284.13333
$ws.Range("I7").Value = 136.2
$ws.Range("K7").Value = 136.2
$ws.Range("M7").Value = -23.19999999999999
$ws.Range("H98").Value = 55000
$ws.Range("J98").Value = 55000
$ws.Range("L98").Value = 55000
$ws.Range("N98").Value = -59492
$ws.Range("H99").Value = 2997.9167
$ws.Range("I99").Value = 2664.111
$ws.Range("J99").Value = 3999.3333
$ws.Range("K99").Value = 2664.111
$ws.Range("L99").Value = 3999.3333
$ws.Range("M99").Value = -1166.111
$ws.Range("N99").Value = -6995.3333
$ws.Range("H105").Value = 958.5
$ws.Range("I105").Value = 869.6667
$ws.Range("K105").Value = 869.6667
$ws.Range("M105").Value = 877.3333
$ws.Range("H126").Value = 2997.9167
$ws.Range("I126").Value = 2664.111
$ws.Range("J126").Value = 3999.3333
$ws.Range("K126").Value = 7992.333
$ws.Range("L126").Value = 11997.9999
$ws.Range("M126").Value = -5522.333
$ws.Range("N126").Value = -16937.9999

$ws = $wb.Sheets("CUL")
$ws.Range("H5").Value = 3230.7778
$ws.Range("I5").Value = 3230.7778
$ws.Range("K5").Value = 9692.3334
$ws.Range("M5").Value = -9580.3334
$ws.Range("H113").Value = 4128.077
$ws.Range("I113").Value = 3962.6667
$ws.Range("J113").Value = 4177.7
$ws.Range("K113").Value = 11888.0001
$ws.Range("L113").Value = 12533.1
$ws.Range("M113").Value = -9718.000100000001
$ws.Range("N113").Value = -16873.1
$ws.Range("H135").Value = 3230.7778
$ws.Range("I135").Value = 3230.7778
$ws.Range("K135").Value = 29077.0002
$ws.Range("M135").Value = -26542.0002

$ws = $wb.Sheets("GSM")
$ws.Range("H122").Value = 3844.1304
$ws.Range("I122").Value = 3774.2632
$ws.Range("K122").Value = 11322.7896
$ws.Range("M122").Value = -8872.7896

$ws = $wb.Sheets("LTW")
$ws.Range("H16").Value = 999.75
$ws.Range("I16").Value = 999.75
$ws.Range("K16").Value = 999.75
$ws.Range("M16").Value = -829.75
$ws.Range("H22").Value = 647.5
$ws.Range("I22").Value = 463.33334
$ws.Range("K22").Value = 463.33334
$ws.Range("M22").Value = -168.33334
$ws.Range("H27").Value = 647.5
$ws.Range("I27").Value = 463.33334
$ws.Range("K27").Value = 463.33334
$ws.Range("M27").Value = -356.33334
$ws.Range("H68").Value = 2598.7778
$ws.Range("I68").Value = 2498
$ws.Range("K68").Value = 2498
$ws.Range("M68").Value = -1749
$ws.Range("H71").Value = 2598.7778
$ws.Range("I71").Value = 2498
$ws.Range("K71").Value = 12490
$ws.Range("M71").Value = -8746
$ws.Range("H82").Value = 2934.375
$ws.Range("J82").Value = 2896.4285
$ws.Range("L82").Value = 2896.4285
$ws.Range("N82").Value = -3618.4285
$ws.Range("H85").Value = 2934.375
$ws.Range("J85").Value = 2896.4285
$ws.Range("L85").Value = 2896.4285
$ws.Range("N85").Value = -5392.4285

$ws = $wb.Sheets("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()
$ws.Range("H107").Value = 1342.8334
$ws.Range("I107").Value = 1681.75
$ws.Range("K107").Value = 5045.25
$ws.Range("M107").Value = -3125.25

